$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl5"
$ws.Range("C2").Value = "Cxcr1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.429319
$ws.Range("H2").Value = 13.287957
$ws.Range("I2").Value = 0.01355902605229267
$ws.Range("J2").Value = 0.01355902605229268
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01218566666666667
$ws.Range("N2").Value = 0.036557
$ws.Range("O2").Value = 0.001498364820294181
$ws.Range("P2").Value = 0.001498364820294181
$ws.Range("Q2").Value = 0.05397420489433333
$ws.Range("R2").Value = 0.4857678440489999
$ws.Range("S2").Value = 0.00002031636763420763
$ws.Range("T2").Value = 0.00002031636763420763

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl5"
$ws.Range("C3").Value = "Cxcr1"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.429319
$ws.Range("H3").Value = 13.287957
$ws.Range("I3").Value = 0.01355902605229267
$ws.Range("J3").Value = 0.01355902605229268
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08128566666666666
$ws.Range("N3").Value = 0.243857
$ws.Range("O3").Value = 0.009994987279658562
$ws.Range("P3").Value = 0.00999498727965856
$ws.Range("Q3").Value = 0.3600401477943332
$ws.Range("R3").Value = 3.240361330149
$ws.Range("S3").Value = 0.0001355222929172243
$ws.Range("T3").Value = 0.0001355222929172243

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cxcl5"
$ws.Range("C4").Value = "Cxcr1"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.429319
$ws.Range("H4").Value = 13.287957
$ws.Range("I4").Value = 0.01355902605229267
$ws.Range("J4").Value = 0.01355902605229268
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.912604999999999
$ws.Range("N4").Value = 23.737815
$ws.Range("O4").Value = 0.9729438112167713
$ws.Range("P4").Value = 0.9729438112167712
$ws.Range("Q4").Value = 35.04745166599499
$ws.Range("R4").Value = 315.427064993955
$ws.Range("S4").Value = 0.01319217048370513
$ws.Range("T4").Value = 0.01319217048370513

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cxcl5"
$ws.Range("C5").Value = "Cxcr1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.429319
$ws.Range("H5").Value = 13.287957
$ws.Range("I5").Value = 0.01355902605229267
$ws.Range("J5").Value = 0.01355902605229268
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.126567
$ws.Range("N5").Value = 0.379701
$ws.Range("O5").Value = 0.015562836683276
$ws.Range("P5").Value = 0.015562836683276
$ws.Range("Q5").Value = 0.5606056178729999
$ws.Range("R5").Value = 5.045450560856999
$ws.Range("S5").Value = 0.0002110169080361154
$ws.Range("T5").Value = 0.0002110169080361154

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cxcl5"
$ws.Range("C6").Value = "Cxcr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 319.3801833333334
$ws.Range("H6").Value = 958.1405500000001
$ws.Range("I6").Value = 0.9776862371851469
$ws.Range("J6").Value = 0.977686237185147
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01218566666666667
$ws.Range("N6").Value = 0.036557
$ws.Range("O6").Value = 0.001498364820294181
$ws.Range("P6").Value = 0.001498364820294181
$ws.Range("Q6").Value = 3.891860454038889
$ws.Range("R6").Value = 35.02674408635
$ws.Range("S6").Value = 0.001464930663084016
$ws.Range("T6").Value = 0.001464930663084017

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cxcl5"
$ws.Range("C7").Value = "Cxcr1"
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 319.3801833333334
$ws.Range("H7").Value = 958.1405500000001
$ws.Range("I7").Value = 0.9776862371851469
$ws.Range("J7").Value = 0.977686237185147
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.08128566666666666
$ws.Range("N7").Value = 0.243857
$ws.Range("O7").Value = 0.009994987279658562
$ws.Range("P7").Value = 0.00999498727965856
$ws.Range("Q7").Value = 25.96103112237222
$ws.Range("R7").Value = 233.64928010135
$ws.Range("S7").Value = 0.009771961504162787
$ws.Range("T7").Value = 0.009771961504162787

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Cxcl5"
$ws.Range("C8").Value = "Cxcr1"
$ws.Range("D8").Value = "Neutrophils"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 319.3801833333334
$ws.Range("H8").Value = 958.1405500000001
$ws.Range("I8").Value = 0.9776862371851469
$ws.Range("J8").Value = 0.977686237185147
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.912604999999999
$ws.Range("N8").Value = 23.737815
$ws.Range("O8").Value = 0.9729438112167713
$ws.Range("P8").Value = 0.9729438112167712
$ws.Range("Q8").Value = 2527.12923554425
$ws.Range("R8").Value = 22744.16311989825
$ws.Range("S8").Value = 0.9512337737811011
$ws.Range("T8").Value = 0.9512337737811011

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Cxcl5"
$ws.Range("C9").Value = "Cxcr1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 319.3801833333334
$ws.Range("H9").Value = 958.1405500000001
$ws.Range("I9").Value = 0.9776862371851469
$ws.Range("J9").Value = 0.977686237185147
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.126567
$ws.Range("N9").Value = 0.379701
$ws.Range("O9").Value = 0.015562836683276
$ws.Range("P9").Value = 0.015562836683276
$ws.Range("Q9").Value = 40.42299166395
$ws.Range("R9").Value = 363.80692497555
$ws.Range("S9").Value = 0.01521557123679908
$ws.Range("T9").Value = 0.01521557123679908

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Cxcl5"
$ws.Range("C10").Value = "Cxcr1"
$ws.Range("D10").Value = "FAPs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.859904666666667
$ws.Range("H10").Value = 8.579714
$ws.Range("I10").Value = 0.008754736762560278
$ws.Range("J10").Value = 0.00875473676256028
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.01218566666666667
$ws.Range("N10").Value = 0.036557
$ws.Range("O10").Value = 0.001498364820294181
$ws.Range("P10").Value = 0.001498364820294181
$ws.Range("Q10").Value = 0.03484984496644444
$ws.Range("R10").Value = 0.313648604698
$ws.Range("S10").Value = 0.00001311778957595649
$ws.Range("T10").Value = 0.00001311778957595649

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Cxcl5"
$ws.Range("C11").Value = "Cxcr1"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.859904666666667
$ws.Range("H11").Value = 8.579714
$ws.Range("I11").Value = 0.008754736762560278
$ws.Range("J11").Value = 0.00875473676256028
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.08128566666666666
$ws.Range("N11").Value = 0.243857
$ws.Range("O11").Value = 0.009994987279658562
$ws.Range("P11").Value = 0.00999498727965856
$ws.Range("Q11").Value = 0.2324692574331111
$ws.Range("R11").Value = 2.092223316898
$ws.Range("S11").Value = 0.00008750348257854916
$ws.Range("T11").Value = 0.00008750348257854916

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Cxcl5"
$ws.Range("C12").Value = "Cxcr1"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.859904666666667
$ws.Range("H12").Value = 8.579714
$ws.Range("I12").Value = 0.008754736762560278
$ws.Range("J12").Value = 0.00875473676256028
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 7.912604999999999
$ws.Range("N12").Value = 23.737815
$ws.Range("O12").Value = 0.9729438112167713
$ws.Range("P12").Value = 0.9729438112167712
$ws.Range("Q12").Value = 22.62929596499
$ws.Range("R12").Value = 203.6636636849099
$ws.Range("S12").Value = 0.008517866951964975
$ws.Range("T12").Value = 0.008517866951964975

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Cxcl5"
$ws.Range("C13").Value = "Cxcr1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.859904666666667
$ws.Range("H13").Value = 8.579714
$ws.Range("I13").Value = 0.008754736762560278
$ws.Range("J13").Value = 0.00875473676256028
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.126567
$ws.Range("N13").Value = 0.379701
$ws.Range("O13").Value = 0.015562836683276
$ws.Range("P13").Value = 0.015562836683276
$ws.Range("Q13").Value = 0.3619695539459999
$ws.Range("R13").Value = 3.257725985514
$ws.Range("S13").Value = 0.0001362485384407981
$ws.Range("T13").Value = 0.0001362485384407981
